$d = $word.ActiveDocument

$d.Content.Find.Execute("Minute 4/10/2010", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Minutes 4/10/2010", 2)
